$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).Cut() | Out-Null
$ws.Columns.Item(3).Insert() | Out-Null
Write-Host "done"
